$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge all data rows (2-10) into a single row (row 2) with combined
# identifiers and a summed forecast value.
$ws.Range("A2").Value = "merged_06GE2U92FA_DTPQ67872X_G8WP29EHC6_HYT3MYM7CY_IMFA21HLV3_JU9OS20S99_PPB56V08LB_R5X6KPETN3_RX38XS00QN"
$ws.Range("B2").Value = "merged_L_M_S_XL_XS_XXL_XXS_onesize"
$ws.Range("C2").Value = "merged_Amsterdam_Berlin_Brussels_Copenhagen_Helsinki_Madrid_Paris_Platform_Rome_Stockholm_Webshop"
$ws.Range("D2").Value = 2024
$ws.Range("E2").Value = 11643.6

# Remove the now-redundant rows (3-10); shifts rows below up automatically.
$ws.Range("A3:E10").EntireRow.Delete()
